# Auto-generated PowerShell COM-interop edit script
# Applies the 'commit species identity from gbif' edits to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v1 = @'
Species mentioned in the text:
1. Wheat (Triticum aestivum)
2. Wheat wild relatives from the Aegilops genus:
   a. Aegilops genus
   b. Aegilops speltoides
   c. Aegilops tauschii
   d. Aegilops sharonensis
   e. Aegilops vavilovii
   f. Aegilops umbellulata
   g. Aegilops peregrina
   h. Aegilops kotschyi
   i. Aegilops longissima
   j. Aegilops cylindrica
   k. Aegilops ukurunduensis
   l. Aegilops speltoides var. ligustica
   m. Aegilops bicornis
Note: The text does not explicitly mention all the species from the Aegilops genus, only 13 accessions were examined.
'@
$ws.Range("F2").Value = $v1
$ws.Range("C3").Value = 'No'
$v2 = @'
- Puccinia triticina (causative agent of leaf rust)
- common winter wheat breeding lines
- isogenic lines carrying genes Lr9, Lr19, Lr22a, Lr22b, and Lr25
- line 99/08-52
- lines 19/06-108, 82/08-43, and 82/08-35
- lines with partial race-specific resistance
- lines with race non-specific resistance
- lines with "slow rusting" type resistance
'@
$ws.Range("F3").Value = $v2
$v3 = @'
Species mentioned in the text:
1. Maize (Zea mays L.)
'@
$ws.Range("F4").Value = $v3
$ws.Range("D5").Value = 'This study compared the growth and yield of ''Muscat Hamburg'' grapes grafted on ''Dog Ridge'' rootstock and self-rooted cuttings in Tamil Nadu, India. The grafted vines showed better growth, yield, and nutrient content in the petiole compared to self-rooted cuttings.'
$v4 = @'
1. Muscat Hamburg grape variety (Vitis vinifera)
2. Dog Ridge rootstock (Vitis champini)
'@
$ws.Range("F5").Value = $v4
$ws.Range("D6").Value = 'Using silicon preparations can promote growth of forest seedlings, particularly oak seedlings affected by oak powdery mildew, and improve their photosynthetic efficiency, biomass, and root nutrition. Spraying with a 2% concentration is most effective.'
$v5 = @'
Species mentioned in the text:
- Pine (Pinus sylvestris)
- Oak (Quercus robur)
'@
$ws.Range("F6").Value = $v5
$ws.Range("D7").Value = 'The study examined the effects of different soil media mixtures with phosphogypsum formulations on the growth of young pine seedlings in Polish forest conditions. The use of phosphogypsum-based preparations did not have toxic effects on the seedlings, and a mixture of phosphogypsum and organic ash showed positive effects on root development. However, longer-term observations are needed to fully understand the impact. Lower dosages of 1 and 2 t/ha appear to be the most promising. Heavy metal testing and monitoring of the microbiome are recommended.'
$v6 = @'
- Pinus sylvestris (Scots pine)
- Pinus silvestris (European pine)
- Organic ash
- Sewage sludge
'@
$ws.Range("F7").Value = $v6
$v7 = @'
1. Guppy (Poecilia reticulata)
2. Artemia salina (brine shrimp)
'@
$ws.Range("F8").Value = $v7
$ws.Range("D9").Value = 'This study explores the growth-survival trade-off in non-phanerophyte species used in dune restoration. Plant traits like leaf dry matter content and floral displays affect this trade-off. Foredune plant species have higher growth but lower survival rates compared to transition dune species. This trade-off can inform cost-effective ecosystem restoration actions.'
$v8 = @'
Species mentioned in the following text:
1. Phanerophyte species
2. Non-phanerophyte species
3. Perennial non-phanerophyte species
4. Plant species of foredune communities
5. Plant species of transition dune communities
'@
$ws.Range("F9").Value = $v8
$ws.Range("C10").Value = 'Unclear'
$ws.Range("D10").ClearContents()
$v9 = @'
Species mentioned in the text:
1. Eelgrass (Zostera marina)
'@
$ws.Range("F10").Value = $v9
$ws.Range("D11").Value = 'The study examines the relationship between seed production and the regeneration of Pinus sibirica, finding that the abundance of seedlings is influenced by the number of Nucifraga birds and vice versa.'
$v10 = @'
Species mentioned in the text:
1. Pinus sibirica Du Tour (Siberian pine)
2. Nucifraga caryocatactes macrorhynchos Brehm C L (thin-billed nutcracker)
'@
$ws.Range("F11").Value = $v10
$ws.Range("C12").Value = 'Yes'
$ws.Range("D12").Value = 'This study examines the distribution of carbon in a seagrass called Posidonia oceanica and how it varies during different seasons and depths in the Mediterranean Sea. The results show that deep-water plants store more carbon during summer and have a different cell wall composition in winter.'
$ws.Range("F12").Value = '1. Posidonia oceanica'
$v11 = @'
1. Alysson spinosus (digger wasp)
2. Hemipteran nymphs or imagines
3. Delphacidae (true hoppers)
4. Cicadellidae (true hopper)
5. lilac plant
6. small-leaved linden plant
7. Metopia argyrocephala (dipteran kleptoparasitic)
8. Sarcophagidae (dipteran) 
9. Alysson melleus (Nearctic digger wasp)
'@
$ws.Range("F13").Value = $v11
$ws.Range("D14").Value = 'This study investigates the occurrence of mosaic diseases caused by Cauliflower mosaic virus (CaMV) and Turnip mosaic virus (TuMV) in cabbage fields in Central Ethiopia. The viruses were found to have a high incidence and wide distribution, affecting cabbage productivity.'
$v12 = @'
- Cabbage (Brassica oleracea var. capitata)
- Cauliflower mosaic virus (CaMV)
- Turnip mosaic virus (TuMV)
- "Habesha gomen" cabbage variety
- "Tikur gomen" cabbage variety
'@
$ws.Range("F14").Value = $v12
$ws.Range("D15").Value = 'The study examined the impact of copper-contaminated soil on barley plants, finding that plant growth depends on copper concentration and soil buffering capacity, with a variable maximum permissible concentration.'
$ws.Range("D16").Value = 'This study tested using hazelnut shell waste as a growth medium for German primroses. Results showed it had a significant impact on appearance and flower weight.'
$ws.Range("E16").Value = 'Yes'
$v13 = @'
1. German primrose (Primula obconica Hance)
2. Sphagnum moss peat (SMP)
3. Hazelnut husk waste (HHW)
'@
$ws.Range("F16").Value = $v13
$ws.Range("D17").Value = 'Soil mixtures derived from gneiss weathering and forest floor types had positive effects on growth, photosynthesis, and nutrient status of chestnut seedlings. Forest floor application is a viable alternative for nutrient management.'
$v14 = @'
In the text, the following species are mentioned:
1. Castanea sativa Mill. (chestnut)
2. Evergreen broad-leaved trees (mentioned in the context of forest floor types)
3. Gneiss (type of soil)
4. Chestnut seedlings
'@
$ws.Range("F17").Value = $v14
$ws.Range("D18").Value = 'This research examines the distribution and resorption efficiency of macroelements in the organs of Rumex alpinus in the Alps and the Giant Mountains. The study found high variability in soil nutrient contents and determined that Rumex alpinus has a high demand for nitrogen, phosphorus, and potassium. However, the plant''s resorption efficiency for these macroelements is lower compared to other plants.'
$v15 = @'
Species mentioned in the text:
1. Rumex alpinus
'@
$ws.Range("F18").Value = $v15
$ws.Range("D19").Value = 'The study found that applying cattle manure and dung beetles to Mediterranean silvopastoral ecosystems increased oak seedling establishment by improving acorn survival and reducing predation.'
$v16 = @'
- Cattle
- Dung beetles
- Invertebrates
- Oaks
- Grassland layer
- Livestock
- Wild predators
- Acorns
- Rodents
'@
$ws.Range("F19").Value = $v16
$ws.Range("C20").Value = 'No'
$ws.Range("D20").ClearContents()
$ws.Range("D21").Value = 'Using a plant defense hormone and a propagation method reduced damage to Norway spruce plants by a bark-feeding insect, showing potential for improved resistance.'
$ws.Range("C22").Value = 'Yes'
$ws.Range("D22").Value = 'A study in Hungary identified a strain of ''Candidatus Phytoplasma asteris'' in Sempervivum plants, causing unusual symptoms, but the plants recovered after a year.'
$v17 = @'
- Sempervivum species
- Cultivars Alpha, Purple Passion, and Silberkarneol
- Onion yellows phytoplasma
- 'Muscari botryoides' aster yellows phytoplasma
'@
$ws.Range("F22").Value = $v17
$ws.Range("D23").Value = 'Neofusicoccum parvum, a fungal species, caused canker and dieback on Alnus glutinosa trees in central Portugal, with a 70% incidence rate. Pathogenicity tests confirmed the first report of its impact on A. glutinosa in Portugal and Europe.'
$v18 = @'
- Neofusicoccum parvum 
- Alnus glutinosa
'@
$ws.Range("F23").Value = $v18
$ws.Range("C24").Value = 'No'
$ws.Range("D24").ClearContents()
$v19 = @'
1. Querciphoma minuta
2. Platanus x hispanica (London plane)
3. Camarosporium complex
'@
$ws.Range("F24").Value = $v19
$ws.Range("D25").Value = 'This paper investigates the quality of Merlot and Blatina grapes and wine in the Trebinje vineyard, demonstrating suitability for producing quality red wines.'
$ws.Range("F25").Value = 'Merlot, Blatina.'
$ws.Range("D26").Value = 'This text discusses the effects of soil contamination with heavy metals from a smelter in Drenas, Kosovo, on maize plants and the potential ecological implications.'
$v20 = @'
- Maize plants (Zea mays)
- Ferronikel smelter
- Fe (Iron)
- Cu (Copper)
- Mn (Manganese)
- Cr (Chromium)
- Cd (Cadmium)
- Ni (Nickel)
- Pb (Lead)
'@
$ws.Range("F26").Value = $v20
$ws.Range("D27").Value = 'Controlled atmosphere temperature treatment (CATT) successfully reduced the number of Phytonemus pallidus (cyclamen mite) in strawberry plants by nearly 100%, offering a potential solution for growers to prevent infestations.'
$v21 = @'
1. Phytonemus pallidus (Cyclamen mite)
2. Fragaria x ananassa Duchesne (Strawberry)
'@
$ws.Range("F27").Value = $v21
$ws.Range("D28").Value = 'This study explores how bumblebees can sense the nutritional status of tomato plants, specifically focusing on the rootstocks and their impact on pollination and yield.'
$ws.Range("E28").Value = 'No'
$ws.Range("F28").ClearContents()
$ws.Range("D29").Value = 'Green compost amendment was found to improve the performance of potato plants grown on Mars regolith simulant as a substrate for cultivation in space.'
$v22 = @'
1. Potato (Solanum tuberosum L., cv. 'Colomba')
2. Mars regolith simulant (MMS-1)
3. Green compost
4. Fluvial sand
5. Red soil from Sicily (RS)
6. Volcanic soil from Campania (VS)
'@
$ws.Range("F29").Value = $v22
$ws.Range("D30").Value = 'Litter leachate affects germination of Norway spruce, but not dwarf pine. Subalpine tall-herb vegetation serves as a filter, impacting woody species seedling composition.'
$v23 = @'
- Norway spruce (Picea abies)
- Dwarf pine (Pinus mugo)
- Calamagrostis villosa
'@
$ws.Range("F30").Value = $v23
$ws.Range("D31").Value = 'Passerine birds in the boreal forest help facilitate seed dispersal and sexual reproduction for Vaccinium shrubs by depositing seeds in tree stumps, providing suitable conditions for seedling establishment.'
$ws.Range("E31").Value = 'Yes'
$v24 = @'
Species mentioned in the text:
1. Vaccinium shrubs (keystone species)
2. Passerine birds (seed dispersers)
3. Bryophytes (required for seedling establishment)
'@
$ws.Range("F31").Value = $v24
$ws.Range("D32").Value = 'Industrial hemp plants have natural resistance to many insect pests due to their physical strength and chemical composition. However, there are still some pests that can cause damage, including cotton bollworm and native budworm. Other sporadic pests include leaf beetles, mirids, and seed-collecting ants. Effective management options include organic, biological, and conventional methods.'
$v25 = @'
- cotton bollworm (Helicoverpa armigera)
- native budworm (H. punctigera)
- Rutherglen bug (Nysius vinitor)
- green vegetable bug (Nezara viridula)
- leaf beetles (Monoleptis australis)
- mirids (Creontiades pallida)
- stem borer (cerambycid beetle)
- seed-collecting ants
- wireworms
- red-legged earth mite (Halotydeus destructor)
- cabbage white butterflies (Pieris rapae)
- miscellaneous grass moths
'@
$ws.Range("F32").Value = $v25
$ws.Range("D33").Value = 'The study measured morphological parameters of Allium ursinum populations in Western Serbia and found that the populations on mountains Povlen and Rudnik have the most productive parameters for agricultural research.'
$ws.Range("F33").Value = 'Allium ursinum'
$ws.Range("D34").Value = 'Field experiments in Slovenia found that nitrogen application levels and cultivars did not interact in their effect on radicchio yield and yield components. The highest fresh weight was achieved with 150 kg N/ha, but dry matter content and firmness of radicchio heads decreased with increased soil N supply. N levels did not significantly affect leaf number.'
$v26 = @'
The different species mentioned in the text are:
1. Nitrogen
2. Radicchio
3. Biotechnical centre
4. Naklo
5. Kranj
6. Slovenia
7. Cultivars: 'Monivip', 'Castel Franco', 'Anivip', 'Foresto', 'Palla rossa', and 'Verona'
8. Headed chicory
9. KAN (calcium ammonium nitrate)
10. Calcium
11. Amonium
12. Fertilizers
13. Crop
'@
$ws.Range("F34").Value = $v26
$ws.Range("C35").Value = 'Unclear'
$v27 = @'
- Sycamore (Acer pseudoplatanus)
- Boxelder maple (Acer negundo)
- Silver maple (Acer saccharinum)
- Sugar maple (Acer saccharum)
- Japanese maple (Acer palmatum)
- Trident maple (Acer buergerianum)
- Paperbark maple (Acer griseum)
- Himalayan maple (Acer oblongum)
'@
$ws.Range("F35").Value = $v27
$ws.Range("D36").Value = 'Breeding assessment of leafiness in Festuca species for forage quality. Evaluation of species, varieties, and ecotypes to select high leafy accessions. Variations in leafiness and the leafiest variety is tall fescue Albena at 59.54%.'
$v28 = @'
Species mentioned in the text:
- Tall fescue
- Meadow fescue
- Red fescue
'@
$ws.Range("F36").Value = $v28
$ws.Range("D37").Value = 'Study focused on determining the phenolic composition of plum fruits grown in Tambov oblast, Russia using high performance liquid chromatography. Plum fruits were found to be a potential source of phenolic compounds.'
$v29 = @'
- Prunus domestica L. (Plum)
- SVG 11-19 (Plum-cherry hybrid)
- 4-Caffeoylquinic acid
- 5-Caffeoylquinic acid
- 3-p-Coumaroylquinic acid
- 3-Caffeoylquinic acid
- Quercitin-3-rutinoside
- Cyanidin-3-glucoside
- Cyanidin-3-rutinoside
- Peonidin-3-glucoside
'@
$ws.Range("F37").Value = $v29
$ws.Range("F38").Value = 'Salvia yangii, Salvia abrotanoides'
$v30 = @'
1. Lemon (Citrus limon)
2. Sour Orange (Citrus aurantium)
3. Yuma Ponderosa lemon (Citrus x jambhiri)
4. Volkameriana (Citrus volkameriana)
'@
$ws.Range("F39").Value = $v30
$v31 = @'
- Apple 
- Bramley's Seedling 
- Falstaff
'@
$ws.Range("F40").Value = $v31
$ws.Range("D41").Value = 'This study examined seed weight and imbibition period of herbaceous peony species native to Serbia, providing preliminary research for future peony germination studies.'
$v32 = @'
- Paeonia tenuifolia (fern leaf peony or steppe peony)
- Paeonia peregrina (Balkan peony or Kosovo peony)
- Paeonia daurica
'@
$ws.Range("F41").Value = $v32
$ws.Range("D42").Value = 'This study evaluated the effects of cultivar and plant origin on minituber production in an aeroponic facility. Agria plants of minituber origin had the highest yield and heaviest tubers.'
$v33 = @'
- Potato (Solanum tuberosum)
- Sinora (Potato cultivar)
- Agria (Potato cultivar)
'@
$ws.Range("F42").Value = $v33
$ws.Range("D43").Value = 'The study found that sodic salinity negatively affects the growth and functioning of the meadow buttercup. The species is best adapted to subsaline habitats and shows limitations at higher salinities. Increases in salinity hampered growth, leaf morphology, and photosynthesis, but not mineral nutrition. Leaf appearance can be a sign of progressing salinity. Ultimately, rising salinity reduces the competitiveness of the species and shifts its strategy to ruderal behavior.'
$v34 = @'
Species mentioned in the text:
1. Meadow buttercup (Ranunculus acris L.)
'@
$ws.Range("F43").Value = $v34
$v35 = @'
1. Coastal Douglas-fir (Pseudotsuga menziesii)
2. Interior Douglas-fir (Pseudotsuga menziesii)
'@
$ws.Range("F44").Value = $v35
$ws.Range("D45").Value = 'The study aimed to assess if pine seedlings could be infested with pinewood nematodes through sawdust. Findings showed that nematodes can move from infested sawdust into damaged seedlings.'
$v36 = @'
- Pine (Pinus sylvestris L.)
- Pinewood nematode (Bursaphelenchus xylophilus Steiner and Buhrer)
- Beetles of the genus Monochamus
'@
$ws.Range("F45").Value = $v36
$ws.Range("D46").Value = 'Various seed treatments, including heat treatment, seed priming, and seed coating, were tested for their potential to improve lucerne seed performance and early field growth. Seed priming with potassium permanganate and chitosan showed the best results for seedling length and emergence dynamics, while coating with bentonite and gypsum had a positive impact on root development. Cinnamon powder improved emergence dynamics, seedling, and shoot length. The combination of priming and coating methods was found to be the most effective in field conditions.'
$ws.Range("E46").Value = 'No'
$ws.Range("F46").ClearContents()
$ws.Range("D47").Value = 'The age of sweet pepper seedlings and the cultivation method have an impact on seedling quality and productivity. 60-day-old seedlings grown by transplanting had higher yield.'
$v37 = @'
Species mentioned in the text:
1. Sweet pepper (Capsicum annuum L.)
'@
$ws.Range("F47").Value = $v37
$ws.Range("C48").Value = 'Yes'
$ws.Range("D48").Value = 'Cd and Pb content in herbs and spices used in Polish cuisine was examined. The results showed low risk of toxicity, but some products exceeded acceptable limits and should be consumed with caution.'
$v38 = @'
species mentioned:
- dried herbs
- fresh herbs
- loose single-component spices
- coriander
- estragon
- watercress
- jiaogulan
- celery
- basil
- dill
'@
$ws.Range("F48").Value = $v38
$ws.Range("C49").Value = 'No'
$ws.Range("D49").ClearContents()
$v39 = @'
- Vitis vinifera cultivars:
  - Black Corinth
  - Carignane
  - Husseine
  - Merlot
  - Muscat Hamburg
  - Palomino
  - Peloursin
  - Primitivo (aka Zinfandel)
  - Thompson Seedless
- Causal fungus:
  - Eutypa lata
'@
$ws.Range("F49").Value = $v39
$v40 = @'
- Quince clonal rootstocks (Province Quince BA29 [BA29], Quince A [QA], and Quince MC [MC])
- European pear (Pyrus communis L.) cultivars: 'Deveci', 'Williams', 'Santa Maria', and 'Abate Fetel'
'@
$ws.Range("F50").Value = $v40
$ws.Range("D51").Value = 'Physiological and biochemical characteristics of drought tolerance in different garden rose genotypes were studied. The best water-holding capabilities were found in cv. ''Borisfen'' and R. hugonis species, while R. indica, R. bracteata, R. rouletti, R. foetida showed instability. Recovery of metabolic processes was observed in R. hugonis, R. bracteata, R. indica, and cv. ''Borisfen'' after mild wilting, but irreversible disturbances occurred in R. gallica, R. indica, and R. bracteata under severe wilting conditions. The highest drought tolerance was found in cv. ''Borisfen'' and R. hugonis.'
$v41 = @'
- Garden roses (genotypes)
- R. hugonis
- R. indica
- R. bracteata
- R. rouletti
- R. foetida
- R. gallica
'@
$ws.Range("F51").Value = $v41
$ws.Range("D52").Value = 'Biodegradable mulching film increases tomato yield and quality. Both types of film tested (biodegradable and polyethylene) showed positive effects, with biodegradable film being more sustainable.'
$ws.Range("E52").Value = 'Yes'
$v42 = @'
1. San Marzano tomato
2. Flavonoids
3. Polyphenols
4. AsA (Ascorbic acid)
'@
$ws.Range("F52").Value = $v42
$ws.Range("D53").Value = 'This study investigated the accumulation and translocation of copper in grapevine grafts grown in different soil media. Both rootstocks showed high copper accumulation, with SO4 rootstock having a higher translocation rate. Soil with high copper content was found to inhibit graft growth.'
$v43 = @'
- Grapevine
- Sauvignon Blanc
- 5BB rootstock
- SO4 rootstock
'@
$ws.Range("F53").Value = $v43
$ws.Range("D54").Value = 'Sea buckthorn is a valuable shrub in Romania, with numerous uses including food, medicine, and forestry. It has high nutritional value and can improve the climate and prevent pollution.'
$ws.Range("F54").Value = 'Hippophae rhamnoides (Sea buckthorn)'
$ws.Range("D55").Value = 'Reciprocal bark exchange between tree species reveals that both bark and wood have different effects on invertebrate communities in mid-decay logs, emphasizing the importance of maintaining deadwood heterogeneity.'
$v44 = @'
1. Araucaria araucana
2. Cryptomeria japonica
3. Picea abies
4. Thuja plicata
5. Chamaecyparis lawsoniana
'@
$ws.Range("F55").Value = $v44
$ws.Range("D56").Value = 'The study compared the chemical composition and antioxidant activity of garlic leaves and bulbs at different stages of development. Garlic leaves contained more protein, fat, ash, vitamin C, and polyphenols than bulbs. Young plants had higher levels of bioactive compounds.'
$ws.Range("F56").Value = '1. Garlic (Allium sativum)'
$ws.Range("C57").Value = 'Yes'
$ws.Range("D57").Value = 'Tomato yellow leaf curl Sardinia virus (TYLCSV) does not transmit through tomato seeds, as genetic material was not detected in seedlings or embryos.'
$ws.Range("F57").Value = 'Tomato yellow leaf curl Sardinia virus, tomato'
$ws.Range("D58").Value = 'Using seed-propagate plug plants of Miscanthus can improve propagation rates and scale up of plantations, while plug design and planting date have significant effects on yield and establishment rates.'
$v45 = @'
Species mentioned in the text:
1. Posidonia oceanica (seagrass)
'@
$ws.Range("F59").Value = $v45
$ws.Range("D60").Value = 'Research on Albizia richardiana plant revealed three phytotoxic compounds that can be used as bioherbicides to control weeds, potentially reducing the need for synthetic chemical herbicides.'
$v46 = @'
1. Albizia richardiana (plant species)
2. Lettuce (weed species)
3. Italian ryegrass (weed species)
4. Lepidium sativum (cress) (weed species)
'@
$ws.Range("F60").Value = $v46
$ws.Range("D61").Value = 'Crop mixtures, specifically those containing faba bean, linseed, or oilseed rape, can increase arthropod biodiversity and abundance, supporting pollination and pest-control ecosystem services without compromising crop yield.'
$v47 = @'
The different species mentioned in the text are:
- Arthropods
- Wheat
- Faba bean
- Linseed
- Oilseed rape
- Grassland biodiversity organisms
- Mass-flowering crops
- Legumes
- Cereals
'@
$ws.Range("F61").Value = $v47
$ws.Range("D62").Value = 'The study examines the recovery and formation of biological soil properties in forest plantations after disturbances such as wildfires using microbiological indicators.'
$ws.Range("F62").Value = 'Pinus sylvestris, Larix sibirica, Ulmus humilis'
$ws.Range("D64").Value = 'Optimizing plant nutrition and cultivation methods can maximize potato yield. The study examines the growth and development of different potato varieties in Forest-Steppe conditions, finding that yield depends on fertilization, planting tuber fraction, and variety characteristics. Maximum yield is achieved with specific fertilization methods and using seed tubers of a certain size. The findings can be applied to improve seed potato production and agribusiness profitability.'
$v48 = @'
- Laperla
- Granada
- Memphis
'@
$ws.Range("F64").Value = $v48
$ws.Range("D65").Value = 'Fungi from Thailand were tested for their ability to biodegrade low-density polyethylene (LDPE) films. Several fungi showed significant growth and degradation of LDPE films, indicating their potential use in plastic degradation.'
$ws.Range("F65").Value = 'Diaporthe italiana, Thyrostroma jaczewskii, Collectotrichum fructicola, Stagonosporopsis citrulli, Aspergillus niger'
$v49 = @'
Citrus x limon (lemon)
Neofusicoccum parvum
Citrus sinensis x Poncirus trifoliata (citrange)
'@
$ws.Range("F66").Value = $v49
$ws.Range("D67").Value = 'The efficacy of Isaria fumosorosea against the box tree moth is low. Host plant phytochemicals may be involved in the moth''s defense against fungal pathogens.'
$v50 = @'
- Isaria fumosorosea
- Cydalima perspectalis
- Buxus sp.
- B. sempervirens
'@
$ws.Range("F67").Value = $v50
$ws.Range("C68").Value = 'Yes'
$ws.Range("D68").Value = 'Defaunation caused by hunting in French Guiana''s rainforests is affecting the functional composition of tree communities, resulting in shifts in leaf and fruit traits and wood density. These changes have long-term ecological consequences.'
$v51 = @'
- Tropical vertebrate populations
- Seed dispersers
- Predators
- Browsers
- Tree recruits
- Tree adults
'@
$ws.Range("F68").Value = $v51
$ws.Range("D69").Value = 'Research explores the relationship between competition and productivity in plant communities, highlighting the confusion between productivity and biomass and the need for further study in this area. The study also examines traits that explain competition and its correlation with diversity.'
$ws.Range("C70").Value = 'Yes'
$ws.Range("D70").Value = 'The study investigates the impact of fish waste on feeding and reproductive ability in Iceland scallops. Results show some differences in feeding but no significant effects on reproductive ability.'
$v52 = @'
Species mentioned in the text:
1. Iceland scallop (Chlamys islandica)
2. Norwegian salmon (Salmo salar)
3. Arctic species
4. Subarctic species
5. Finfish (unspecified species)
'@
$ws.Range("F70").Value = $v52
$ws.Range("D71").Value = 'This study compared the chemical composition of Syringa vulgaris (lilac) and soil in abandoned cemeteries. The plant''s tissue composition was more influenced by soil substrate and formation process than burial sites.'
$ws.Range("F71").Value = 'Syringa vulgaris (lilac), Ca, Na, Mg, Al, Fe, Zn, Cd, Pb'
$ws.Range("D72").Value = 'Real-time PCR was used to detect and quantify oomycetes in ornamental plants bought from various sources. Oomycete DNA was found in composts, roots, and filters, with the highest quantities detected using the ITS probe. No differences were found between plants purchased online or from traditional retailers.'
$ws.Range("D73").Value = 'Study investigates the effects of warmer temperatures on the transcriptomes of different strawberry ecotypes during asexual and sexual reproduction, highlighting differences in gene expression and splicing isoforms.'
$v53 = @'
1. Fragaria vesca (strawberry)
2. SOC1
3. LHY
4. SVP
'@
$ws.Range("F73").Value = $v53
$ws.Range("D74").Value = 'Changes in volatile profiles of arborvitae from drought and insect infestation were analyzed. Key components specific to stressed trees were identified, which could be useful for pest management.'
$v54 = @'
1. Arborvitae (Thuja occidentalis)
2. Cypress bark beetle (Phloeosinus aubei)
3. Coleoptera
4. Curculionidae
5. Scolytinae
6. Cupressaceae
7. Pinene (a-pinene)
8. Thujene (a-thujene)
9. Thujone (a-thujone)
10. Beta-pinene
11. Myrcene
12. Limonene
13. P-cymene
14. Camphene
15. Fenchene
16. Frass
17. Fenchone
'@
$ws.Range("F74").Value = $v54
$ws.Range("F75").Value = 'Salvia fruticosa, Malva sylvestris, Taraxacum officinale, Plantago ovata, Tanacetum parthenium, Allium sativum'
$ws.Range("C76").Value = 'Yes'
$ws.Range("D76").Value = 'First report of canna yellow streak virus in Iran''s canna plants causing severe symptoms. The virus is similar to isolates from Russia and UK.'
$ws.Range("F76").Value = 'Canna indica, Potyvirus, canna yellow streak virus (CaYSV), Cannaceae family.'
$ws.Range("D77").Value = 'Drought and salt stress affect the emission rates and composition of biogenic volatile organic compounds (BVOCs) in urban trees, with specific changes in BVOC blends and herbivore-related bouquets.'
$ws.Range("F77").Value = 'Quercus robur (oak), Fagus sylvatica (beech), Betula pendula (silver birch), Carpinus betulus (hornbeam)'
$ws.Range("D78").Value = 'Severe dieback and mortality of wild olive trees in Italy caused by Phytophthora species, including the first report of P. bilorbang on wild olive trees.'
$v55 = @'
1. Wild olive trees (Olea europaea var. sylvestris)
2. Ceratonia siliqua (used as bait)
3. Phytophthora species
4. P. bilorbang (CBS131653)
5. P. pseudocryptogea (CBS139749)
'@
$ws.Range("F78").Value = $v55
$ws.Range("D79").Value = 'The study examined the effects of saltwater irrigation on Juglans regia ''Sorrento'' plants. The plants showed no significant morphological or growth changes, indicating potential for cultivation in salinized environments.'
$ws.Range("D80").Value = 'Masting increases fruit production in rowan trees, leading to higher seedling recruitment both near and far from rowans, supporting predator satiation and animal dispersal hypotheses.'
$v56 = @'
Species mentioned in the text:
1. Rowan (Sorbus aucuparia) - fleshy-fruited tree
2. Heterospecifics - refers to non-conspecific trees in the vicinity
3. Frugivores - animals that eat fruits and disperse seeds
'@
$ws.Range("F80").Value = $v56
$ws.Range("D81").Value = 'This article discusses the successful reproduction of four species of the Elaeagnus genus using soft cuttings and growth regulators, which could increase biodiversity in the genus.'
$ws.Range("F81").Value = 'Four species of the genus Elaeagnus L.'
$ws.Range("C82").Value = 'No'
$ws.Range("D82").ClearContents()
$v57 = @'
The species mentioned in the text are:
1. Bryophytes (group of organisms)
2. European Red List of bryophytes (specific list)
3. Biological traits (related to life history, growth habit, sexual and vegetative reproduction)
4. Ecological traits (indicator values, substrate, and habitat)
5. Bioclimatic variables (based on the species' European range)
'@
$ws.Range("F82").Value = $v57
$ws.Range("D83").Value = 'A study on Tarenaya cultivars found low morphological diversity but clear genetic differentiation between seed-propagated and vegetatively propagated cultivars. Vegetatively propagated cultivars showed better performance.'
$v58 = @'
Species mentioned in the text:
1. Tarenaya hassleriana (syn. Cleome spinosa)
2. Tarenaya boliviensis
'@
$ws.Range("F83").Value = $v58
$ws.Range("D84").Value = 'Interspecific hybridization between tulip cultivars and T. altaica showed varying compatibility in terms of germination, fruit-setting, and seed formation. T. altaica had the highest germination rate compared to other parents. Successful crosses included ''Heart of Poland'' x T. altaica and ''Bolroyal Dream'' x T. altaica.'
$v59 = @'
1. Tulipa altaica 
2. 'Heart of Poland' (tulip cultivar)
3. 'Golden Parade' (tulip cultivar)
4. 'Purple Dream' (tulip cultivar)
5. 'Crystal Star' (tulip cultivar)
6. 'Bolroyal Dream' (tulip cultivar)
'@
$ws.Range("F84").Value = $v59
$ws.Range("D86").Value = 'A greenhouse experiment found that increasing drought stress reduced plant growth, yield, and quality in two Sicilian rose species. Identifying drought-tolerant roses could benefit nursery production in water-scarce regions.'
$ws.Range("F86").Value = 'Rosa canina and Rosa sempervirens'
$ws.Range("C87").Value = 'Yes'
$ws.Range("D87").Value = 'An ecological study in Istanbul evaluates how the physical properties of plant leaves affect sound absorption and attenuation in urban areas, using experimental data and statistical analysis.'
$v60 = @'
Species mentioned in the text:
1. English Ivy (Hedera helix L.)
2. Horse Chestnut (Aesculus hippocastanum L.)
3. Hortensia (Hydrangea macrophylla (Thunb.) Ser.)
4. Japanese Privet (Ligustrum japonicum Thunb.)
5. Laurel (Laurus nobilis L.)
6. Linden (Tilia tomentosa Moench)
7. Magnolia (Magnolia grandiflora L.)
8. Osmanthus (Osmanthus heterophyllus (G. Don) P.S. Green)
9. Plane Tree (Platanus orientalis L.)
10. Cherry Laurel (Prunus laurocerasus L.)
'@
$ws.Range("F87").Value = $v60
$ws.Range("D88").Value = 'The study explored the effects of species and varietal differences on pollination and hybrid seedling survival between apple and pear, highlighting prezygotic and postzygotic barriers, and providing potential markers for intergeneric hybrids.'
$v61 = @'
- Apple (species)
- Pear (species)
- European pear (variety)
- Japanese pear (variety)
- Chinese pear (variety)
'@
$ws.Range("F88").Value = $v61
$ws.Range("D89").Value = 'Fish effluents were tested as fertilisers for onion crops and found to enhance soil fertility and yield. The use of fish effluents may challenge current regulations on organic agriculture.'
$v62 = @'
- Onion (Allium cepa)
- Flavobacterium
- Pseudarthrobacter
- Sphingomonas
- Massilia
- Nitrososphaera
- Pseudomonas
- Nocardioides
'@
$ws.Range("F89").Value = $v62
$ws.Range("D90").Value = 'The study examines the growth and biomass production of Norway spruce at lower altitudes in the Czech Republic, and recommends specific conditions for cultivation.'
$ws.Range("F90").Value = 'Norway spruce (Picea abies L. Karst)'
$ws.Range("D91").Value = 'Water limitation during reproduction of Frangula alnus shrubs had transgenerational effects, increasing germination and advancing seedling emergence, with population differentiation in timing and germination stability.'
$ws.Range("F91").Value = '1. Frangula alnus'
$ws.Range("D92").Value = 'Study in Portugal assessed the herbicidal potential of waste extracts from invasive plants (Acacia dealbata bark, Oxalis pes-caprae biomass) and agri-food (spent coffee grounds) on common urban weeds. Extracts reduced pre-emergence performance of weeds but had limited post-emergence effect. Soil neutralized pre-emergence effect, suggesting the need for additional synthetic herbicides in areas without soil. Combining bioherbicides and commercial formulations promotes sustainability.'
$ws.Range("D93").Value = 'Glasshouse experiments showed that fertilization with bioslurry at 100% and 200% rates improved the growth and quality of Swiss chard seedlings compared to inorganic fertilizers.'
$v63 = @'
- Swiss chard (Beta vulgaris L.)
- Star 1801 (cultivar of Swiss chard)
- Fordhook giant (cultivar of Swiss chard)
'@
$ws.Range("F93").Value = $v63
$ws.Range("D94").Value = 'Research in Mexico has identified Tomato brown rugose fruit virus (ToBRFV) in greenhouses and found that it can be transmitted through seed coats, seedlings, and water sources. Chemical and heat treatments were found to be effective in sanitizing growing media and utensils.'
$v64 = @'
1. Tomato brown rugose fruit virus (ToBRFV)
2. Nicotiana rustica
'@
$ws.Range("F94").Value = $v64
$ws.Range("D95").Value = 'The Canadian Plant Disease Survey found that in 2022, cold, wet soils caused nutrient deficiency in plants, while cool, rainy weather led to high levels of botrytis in berries. Hot, dry weather in the following months caused heat and drought stress. Powdery mildew was also observed on many plants. No new diseases were found.'
$v65 = @'
Species mentioned in the text:
1. Berry plants 
2. Vegetable plants 
3. Ornamental nursery plants 
4. Landscape plants 
5. Botrytis 
6. Powdery mildew 
7. Big-leaf maple (Acer macrophyllum)
'@
$ws.Range("F95").Value = $v65
$v66 = @'
1. Black turmeric (Curcuma caesia)
2. Staphylococcus aureus (bacterial strain)
3. Escherichia coli (bacterial strain)
4. Candida albicans (fungus strain)
5. Human peroxiredoxin 5 (protein)
6. Tyrosyl-tRNA synthetase from Staphylococcus aureus (protein)
7. Glucosamine 6-phosphate synthase from Escherichia coli (protein)
8. Zingiberaceae family
'@
$ws.Range("F96").Value = $v66
$ws.Range("D97").Value = 'The study evaluated the long-term effects of living mulches on apple tree growth, yield, and fruit quality. It found that the presence of living mulch affected young tree growth and reduced fruit yield, but improved fruit coloration. The use of semidwarf rootstock, delayed sowing of the cover crop, and high-quality nursery stock were important for apple tree tolerance to living mulch.'
$v67 = @'
Species mentioned in the text:
1. Colonial bent grass (Agrostis vulgaris)
2. White clover (Trifolium repens)
3. Blue fescue (Festuca ovina)
4. Apple tree 'Ligol' (Malus)
5. M.26 rootstock
6. M.9 rootstock
7. P 60 rootstock
8. P 2 rootstock
9. P 16 rootstock
10. P 22 rootstock
'@
$ws.Range("F97").Value = $v67
$ws.Range("D98").Value = 'Study investigates how different pear rootstocks and cultivars affect leaf chlorophyll content, growth, and yield. Results show significant differences and suggest canopy management can improve chlorophyll content and yield.'
$ws.Range("E98").Value = 'No'
$ws.Range("F98").ClearContents()
$ws.Range("D99").Value = 'The study investigated the impact of different nitrogen fertilizers on lettuce yield and quality. Ammonium nitrate resulted in the highest average weight and urea had the highest vitamin C content. Red lettuce varieties had lower vitamin C content and higher nitrate accumulation.'
$v68 = @'
Lettuce (Lactuca sativa L.)
Varieties of lettuce: unidentified
Nitrogen fertilizers: ammonium sulphate, ammonium nitrate, urea
'@
$ws.Range("F99").Value = $v68
$ws.Range("D100").Value = 'Nitrogen fertilization in apple trees can boost plant growth but also increase the risk of European canker disease. Research is being conducted to determine the cause and find management strategies.'
$v69 = @'
1. Apple tree (Malus domestica)
2. European canker (Neonectria ditissima)
3. Fungal pathogen (Neonectria ditissima)
4. Urea (a nitrogen-containing product)
5. cv. Gala trees (a variety of apple tree)
6. PDA (Potato Dextrose Agar, a growth medium)
7. Leaf buds
8. Spores
9. Soil
'@
$ws.Range("F100").Value = $v69
$ws.Range("F101").Value = 'Neonectria ditissima'

Write-Host "edit complete"
